$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "JOHNNY CASH"
$ws.Range("A3").Value = $null
$ws.Range("A4").Value = $null
$ws.Range("A5").Value = $null
$ws.Range("A8").Select()
